$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 187, shifting existing rows 187:254 down to 188:255
$ws.Rows("187:187").Insert()

# Populate the newly inserted row 187 with the new data record
$ws.Cells.Item(187, 1).Value  = 10
$ws.Cells.Item(187, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(187, 3).Value  = "La Araucanía"
$ws.Cells.Item(187, 4).Value  = 45027
$ws.Cells.Item(187, 5).Value  = 9
$ws.Cells.Item(187, 6).Value  = "Fruta"
$ws.Cells.Item(187, 7).Value  = 100104
$ws.Cells.Item(187, 8).Value  = "Frutos de pepita"
$ws.Cells.Item(187, 9).Value  = 100104003
$ws.Cells.Item(187, 10).Value = "Membrillo"
$ws.Cells.Item(187, 11).Value = "Champion"
$ws.Cells.Item(187, 12).Value = "Primera"
$ws.Cells.Item(187, 13).Value = 150
$ws.Cells.Item(187, 14).Value = 15000
$ws.Cells.Item(187, 15).Value = 15000
$ws.Cells.Item(187, 16).Value = 15000
$ws.Cells.Item(187, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item(187, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(187, 19).Value = 833
$ws.Cells.Item(187, 20).Value = 18
